$d = $word.ActiveDocument

# --- Step 1: locate the "14 [POS] Front End User Registration :" paragraph
#     and append a trailing space run to it (" " -> "14 [POS] ... : ") ---
$searchRange = $d.Content
$found = $searchRange.Find.Execute("14 [POS] Front End User Registration :", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target paragraph '14 [POS] Front End User Registration :'"
}
$targetPara = $searchRange.Paragraphs(1)
$paraRange = $targetPara.Range
# exclude the trailing paragraph mark so the new run lands inside the paragraph
$textOnlyRange = $d.Range($paraRange.Start, $paraRange.End - 1)
$textOnlyRange.InsertAfter(" ")

# --- Step 2: insert the new paragraphs (bulleted list items + the new
#     "15 [POS] Front End Send OTP" heading paragraph) right after it,
#     just before the document's final (empty) paragraph ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$snippetXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>registration form : registration-form.blade.php</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>set id in input fields</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">insert </w:t></w:r><w:r><w:t>&lt;form action="" onsubmit="return false"&gt;</w:t></w:r><w:r><w:t xml:space="preserve"> to restrict page reloading</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>do necessary coding</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>&amp; test</w:t></w:r></w:p>
<w:p><w:r><w:t>15 [POS] Front End Send OTP</w:t></w:r><w:r><w:t xml:space="preserve"> : for Password Reset  </w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>xx</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertionPoint.InsertXML($snippetXml) | Out-Null
